$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that must remain literal text
# (matches the source inlineStr cells, e.g. "1.000", "29.386.89").
# Force text format on D before assigning so Excel does not coerce
# these into actual numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.386.89'
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.847.51'
$ws.Range("E3").Value = '  -0.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9987'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.45'
$ws.Range("E5").Value = '  -0.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6309'
$ws.Range("E6").Value = '  +0.72%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07555'
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2961'
$ws.Range("E9").Value = '  -0.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.48'
$ws.Range("E10").Value = '  +0.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07716'
$ws.Range("E11").Value = '  +0.48%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.871.98'
$ws.Range("E12").Value = '  -0.69%  '
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6852'
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001004'
$ws.Range("E15").Value = '  +3.12%  '
$ws.Range("E16").Value = '  -0.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.153'
$ws.Range("E17").Value = '  -1.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.432.10'
$ws.Range("E18").Value = '  -0.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '229.73'
$ws.Range("E19").Value = '  -2.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.44'
$ws.Range("E20").Value = '  -0.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.555'
$ws.Range("E22").Value = '  -0.88%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '157.01'
$ws.Range("E24").Value = '  +0.75%  '
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.379'
$ws.Range("E26").Value = '  -0.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.67'
$ws.Range("E27").Value = '  -0.26%  '
$ws.Range("E28").Value = '  -0.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05725'
$ws.Range("E29").Value = '  -1.98%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.255'
$ws.Range("E30").Value = '  -0.35%  '
$ws.Range("E31").Value = '  +0.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.024'
$ws.Range("E32").Value = '  -0.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.851'
$ws.Range("E33").Value = '  -1.92%  '
$ws.Range("E34").Value = '  -1.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7163'
$ws.Range("E35").Value = '  -0.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.589'
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.251.19'
$ws.Range("E37").Value = '  +1.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01817'
$ws.Range("E38").Value = '  +2.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.781'
$ws.Range("E39").Value = '  -0.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9126'
$ws.Range("E40").Value = '  +0.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.185'
$ws.Range("E41").Value = '  +0.92%  '
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.997.71'
$ws.Range("E43").Value = '  -2.93%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.78'
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '66.25'
$ws.Range("E45").Value = '  -1.63%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.076'
$ws.Range("E46").Value = '  -3.26%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.149'
$ws.Range("E47").Value = '  +0.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4029'
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000117'
$ws.Range("E49").Value = '  -0.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.702'
$ws.Range("E50").Value = '  -0.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1132'
$ws.Range("E51").Value = '  +1.17%  '
